$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "smith"
$ws.Range("B4").Value = "com"
$ws.Range("C4").Value = "skdjfskdjfsdkfjds"
$ws.Range("D4").Value = "dfsjdfosijdfsl sdfjsdkj sdfjl"

$ws.Range("B8").Select()
